# ACDC2019-Naturalist-SpecsAnalyses.xlsx
# "AutoDS Enrichissement E/S données + cohérence jeux de tests"
#
# 1) Add the missing accents to the French header labels used across the
#    sample/model/params sheets (Espece -> Espèce, Duree -> Durée,
#    FonctionCle -> FonctionClé, SerieAjust -> SérieAjust). The header
#    positions/meaning are unchanged, only the spelling gains accents.
# 2) Refresh each sheet's remembered selection / active cell, and move the
#    active workbook tab from "Params2_expl" back to "Echant1_impl".

$wb = $excel.ActiveWorkbook

$wsEchant1 = $wb.Worksheets.Item("Echant1_impl")
$wsEchant2 = $wb.Worksheets.Item("Echant2_impl")
$wsModl    = $wb.Worksheets.Item("Modl_impl")
$wsParams1 = $wb.Worksheets.Item("Params1_expl")
$wsParams2 = $wb.Worksheets.Item("Params2_expl")

# --- Header label spelling fixes (add French accents) ---------------------

$wsEchant1.Range("A1").Value = "Espèce"
$wsEchant1.Range("D1").Value = "Durée"

$wsEchant2.Range("A1").Value = "Espèce"
$wsEchant2.Range("D1").Value = "Durée"

$wsModl.Range("A1").Value = "FonctionClé"
$wsModl.Range("B1").Value = "SérieAjust"

$wsParams1.Range("A1").Value = "Espèce"
$wsParams1.Range("D1").Value = "Durée"

$wsParams2.Range("A1").Value = "Espèce"
$wsParams2.Range("D1").Value = "Durée"

# --- Selections / active sheet ---------------------------------------------
# Leave each non-final sheet with its new remembered selection, then select
# Echant1_impl last so it ends up as the active tab, matching the diff.

$wsParams2.Range("A2").Select()
$wsParams1.Range("H13").Select()
$wsModl.Range("F10").Select()
$wsEchant2.Range("E9").Select()
$wsEchant1.Range("C9").Select()
